# Apply "Add data for 2023-04-11" update to cta-violent-crime-ytd.xlsx
# 1) Reorder tabs so "Irving Park" precedes "Galewood"
# 2) Expand/update the "Irving Park" sheet with a new 2016 data column + revised counts
# 3) Apply incremented totals across several neighborhood sheets + the two summary sheets

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Move "Irving Park" so it sits right before "Galewood" in the tab order
# ---------------------------------------------------------------------------
$galewood = $wb.Worksheets.Item("Galewood")
$irving = $wb.Worksheets.Item("Irving Park")
$irving.Move($galewood)

# ---------------------------------------------------------------------------
# 2. Rebuild "Irving Park" sheet data (new 2016 column inserted before 2017,
#    shifting everything right by one column, plus a few updated counts)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Irving Park")

# Give the brand-new header cell H1 the same style as the rest of row 1
# (bold/bordered/centered) by copying an existing header cell onto it.
$ws.Range("G1").Copy($ws.Range("H1"))

# Row 1 - year headers (shifted right by one, 2016 inserted at B1)
$ws.Range("B1").Value = 2016
$ws.Range("C1").Value = 2017
$ws.Range("D1").Value = 2018
$ws.Range("E1").Value = 2020
$ws.Range("F1").Value = 2021
$ws.Range("G1").Value = 2022
$ws.Range("H1").Value = 2023

# Row 2 - Aggravated Assault
$ws.Range("B2").Value = 1
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 1
$ws.Range("F2").ClearContents()
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 1

# Row 3 - Robbery
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 1
$ws.Range("D3").ClearContents()
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("G3").ClearContents()
$ws.Range("H3").Value = 1

# Row 4 - Total
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 2

# ---------------------------------------------------------------------------
# 3. Updated counts on other neighborhood sheets + the two summary sheets
# ---------------------------------------------------------------------------

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("C2").Value = 12
$ws.Range("D2").Value = 19
$ws.Range("H3").Value = 24
$ws.Range("I3").Value = 44
$ws.Range("J3").Value = 51
$ws.Range("D6").Value = 119
$ws.Range("H6").Value = 107
$ws.Range("C7").Value = 159
$ws.Range("D7").Value = 180
$ws.Range("H7").Value = 165
$ws.Range("I7").Value = 216
$ws.Range("J7").Value = 206

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I3").Value = 3
$ws.Range("I6").Value = 13

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("D2").Value = 1
$ws.Range("D6").Value = 7

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I3").Value = 3
$ws.Range("I5").Value = 7

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J3").Value = 1
$ws.Range("J5").Value = 3

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("H6").Value = 16
$ws.Range("H7").Value = 21

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("H3").Value = 1
$ws.Range("H5").Value = 7

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("C4").Value = 2
$ws.Range("C5").Value = 4

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("D26").Value = 4
$ws.Range("I29").Value = 13
$ws.Range("D33").Value = 7
$ws.Range("J38").Value = 7
$ws.Range("C40").Value = 1
$ws.Range("H50").Value = 21
$ws.Range("H58").Value = 1
$ws.Range("H62").Value = 7
$ws.Range("J77").Value = 3
$ws.Range("C92").Value = 159
$ws.Range("D92").Value = 180
$ws.Range("H92").Value = 165
$ws.Range("I92").Value = 216
$ws.Range("J92").Value = 206

Write-Output "done"
